$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: remove the "is_question_in_expected" boolean column (old column E) ---
# This shifts the old F (Expected Answer) / G (similarity_score) columns left to E / F.
$ws.Columns.Item(5).Delete()

# After the delete, the sheet has columns:
#   A = row index, B = Question, C = Model Name, D = Response,
#   E = Expected Answer, F = similarity_score
# and 2 data rows (old rows 2 & 3 - the "How many curves..." Q&A for the two models).

# --- Step 2: relocate the existing "curves" rows down to their final positions (rows 4 & 5) ---
# Cut+Paste (rather than plain value copy) so the formatting (e.g. the bordered/bold
# style on column A) travels with the cells.
$ws.Range("A2:F2").Cut($ws.Range("A4:F4"))
$ws.Range("A3:F3").Cut($ws.Range("A5:F5"))

# --- Step 3: row 2 - "What are the formats for loading a text file?" (deepseek1.5) ---
$ws.Cells.Item(2,1).Value2 = 0
$ws.Cells.Item(2,2).Value2 = "What are the formats for loading a text file?"
$ws.Cells.Item(2,3).Value2 = "deepseek1.5"
$ws.Cells.Item(2,4).Value2 = "GEO can load several different types of ASCII files, such as CWLAS (LAS), tab-delimited, space-delimited and comma-delimited. The preferred data formats to request from your vendors are wholly structured LAS format or compatible-XML format."
$ws.Cells.Item(2,6).Value2 = -0.01678333058953285

# --- Step 4: row 3 - "How many tracks can you define in one ODF...?" (llama3.2:latest) ---
$ws.Cells.Item(3,1).Value2 = 1
$ws.Cells.Item(3,2).Value2 = "How many tracks can you define in one ODF (Output Database Type File)?"
$ws.Cells.Item(3,3).Value2 = "llama3.2:latest"
$ws.Cells.Item(3,4).Value2 = "According to the GEO application documentation, there is no specific limit mentioned on the number of tracks that can be defined in one ODF. However, it's recommended to keep the track count reasonable for efficient data management and analysis."
$ws.Cells.Item(3,6).Value2 = 0.005372118670493364

# --- Step 5: fix up the row index / model name on the relocated "curves" rows (now rows 4 & 5) ---
$ws.Cells.Item(4,1).Value2 = 2
$ws.Cells.Item(4,3).Value2 = "deepseek1.5"

$ws.Cells.Item(5,1).Value2 = 3
$ws.Cells.Item(5,3).Value2 = "llama3.2:latest"

# --- Step 6: row 6 - "What's the maximum number of characters in a single text entry?" ---
# This is a brand new row, so first copy the column-A style (bordered/bold) from an
# existing data row so the new index cell matches the sheet's formatting convention.
$ws.Cells.Item(4,1).Copy()
$ws.Cells.Item(6,1).PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Cells.Item(6,1).Value2 = 4
$ws.Cells.Item(6,2).Value2 = "What's the maximum number of characters in a single text entry?"
$ws.Cells.Item(6,3).Value2 = "llama3.2:latest"
$ws.Cells.Item(6,4).Value2 = "The maximum number of characters in a single text entry is 32000."
$ws.Cells.Item(6,6).Value2 = 0.04994607716798782
